$wb = $excel.ActiveWorkbook

# Sheet "CaseDetailStat": F2 size value -> formatted "105.75 KB" (was raw "105.751953125")
$wsDetail = $wb.Worksheets.Item("CaseDetailStat")
$wsDetail.Range("F2").Value = "105.75 KB"

# Sheet "CaseDetailStat_Message": A28 holds the Cypher query text used to
# produce the CaseDetailStat sheet -> update to the new rounded/unit-formatted query
$wsMsg = $wb.Worksheets.Item("CaseDetailStat_Message")
$wsMsg.Range("A28").Value = 'MATCH (f:file)-[*]->(c:case) WITH DISTINCT(f) AS f, c MATCH (f)-->(parent)
WHERE c.case_id IN [''NCATS-COP01CCB050022'']
WITH
[''Bytes'', ''KB'', ''MB'', ''GB'', ''TB''] AS units,
toInteger(floor(log(f.file_size)/log(1024))) as i,
2 as precision,
f,parent
WITH f.file_size /(1024^i) AS value, 10^precision AS factor, units[i] as unit,f,parent
RETURN f.file_name AS `File Name` ,f.file_type AS `File Type`,head(labels(parent)) AS `Association`, f.file_description AS `Description`,f.file_format AS Format,round(factor * value)/factor+ +unit AS Size'
